# Apply updated crypto price/volume figures (data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.993.76"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.338.81"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.32"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  -3.95%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.74"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.45"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.81%  "
$ws.Range("D16").Value = "2.349.50"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.831"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "42.924.11"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "0.0₃0910"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.60"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.93"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.47"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.32"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0723"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("D44").Value = "2.020.58"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.68"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.71"
$ws.Range("D49").ClearFormats()
$ws.Range("D51").Value = "2.563.87"
$ws.Range("E51").Value = "  +1.18%  "
